$d = $word.ActiveDocument

function Get-WdColor($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# --- New paragraph style "Abstract Title" (styleId AbstractTitle), based on
#     Normal, followed by Abstract. ---
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1         # wdAlignParagraphCenter
$abstractTitle.ParagraphFormat.SpaceBefore = 15      # 300 twips
$abstractTitle.ParagraphFormat.SpaceAfter = 0        # 0 twips

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = Get-WdColor 0x34 0x5A 0x8A   # 345A8A

# --- Abstract style: reduce space-before from 15pt (300 twips) to 5pt (100 twips) ---
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- ImportTok character style: add green, bold run formatting ---
$importTok = $d.Styles("ImportTok")
$importTok.Font.Color = Get-WdColor 0x00 0x80 0x00        # 008000
$importTok.Font.Bold = $true

# --- BuiltInTok character style: add green run formatting ---
$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = Get-WdColor 0x00 0x80 0x00        # 008000
